$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.005", "24.017.07").
# Force text format before assignment so Excel keeps them as literal strings
# instead of auto-converting to numbers/dates, then restore the default "Normal"
# style so no formatting changes leak into the saved cell (matches source diff,
# which only changes cell values/text, not styles).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.017.07'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.95%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.641.61'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.70%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.002'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '306.96'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.86%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3879'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.97%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3840'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.88%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.001'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.344'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -8.72%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '48.92'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -7.96%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08432'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -4.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '23.70'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -9.07%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.097'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.00001276'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -5.59%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.451'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -6.47%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.638.17'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -5.10%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '94.43'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06953'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.61'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.897'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.74%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.003'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '13.57'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.020.37'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.96%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.323'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.58%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.681'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -8.93%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.46'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.73%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '157.80'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.626'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '141.21'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -7.19%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.230'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -13.91%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.453'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -7.82%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.818.11'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.53%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.993'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07986'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.26%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02899'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -8.35%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9582'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -8.84%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2683'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -7.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.09196'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.461'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.904'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -9.83%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7556'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -8.52%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.02'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -7.17%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.00'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -6.24%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6866'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.04%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.467'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -8.19%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.081'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.05%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.08337'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -9.96%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '133.05'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.87%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.253'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -10.80%  '
